# Atualização de bases das ligas, do dia: 19-04-2024 às 00:38
# Swap the full data rows (columns B:AC) for the following row pairs. The
# running index in column A stays put; every other field (id, teams,
# odds, etc.) for these matches moves to the other row of the pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(5, 7),
    @(20, 21),
    @(67, 68),
    @(120, 121),
    @(125, 126)
)

foreach ($pair in $pairs) {
    $row1 = $pair[0]
    $row2 = $pair[1]

    $range1 = $ws.Range("B$($row1):AC$($row1)")
    $range2 = $ws.Range("B$($row2):AC$($row2)")

    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value = $values2
    $range2.Value = $values1
}
